$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New text/number values for rows 8-11 ----
$A8 = 'SCRIPT/G01P03A/um1312.ssb'
$B8 = 319
$C8 = ' A bad Pokémon like [CS:N]Grovyle[CR]\ncan\''t be allowed to roam free.'
$D8 = ' Таким Покемонам-негодяям как\n[CS:N]Гровайл[CR] нельзя позволить бродить по миру\nбезнаказанно.'
$E8 = ' Óàëéí Ðïìåíïîàí-îåãïäÿÿí ëàë\n[CS:N]Ãñïâàêì[CR] îåìûèÿ ðïèâïìéóû áñïäéóû ðï íéñô\náåèîàëàèàîîï.'
$A9 = 'SCRIPT/G01P03A/um1404.ssb'
$B9 = 322
$C9 = ' The capture of [CS:N]Grovyle[CR]...[K]\nWe\''re here to help!'
$D9 = ' Мы пришли, чтобы помочь вам...[K]\nСхватить [CS:N]Гровайла[CR]!'
$E9 = ' Íú ðñéšìé, œóïáú ðïíïœû âàí...[K]\nÒöâàóéóû [CS:N]Ãñïâàêìà[CR]!'
$A10 = 'SCRIPT/T01P02A/um1603.ssb'
$B10 = 297
$C10 = ' We heard…'
$D10 = ' Мы узнали...'
$E10 = ' Íú ôèîàìé...'
$B11 = 300
$C11 = ' Yes, we heard why [CS:N]Grovyle[CR] was\nstealing Time Gears in the first place.'
$D11 = ' Да, мы узнали зачем [CS:N]Гровайл[CR]\nпохищал Шестерни Времени.'
$E11 = ' Äà, íú ôèîàìé èàœåí [CS:N]Ãñïâàêì[CR]\nðïöéþàì Šåòóåñîé Âñåíåîé.'


# ---- Copy formatting (and placeholder values, overwritten below) from existing rows ----
# Row 8 uses the same banding style as row 6 (no border)
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A8:E8").PasteSpecial(-4104) | Out-Null

# Row 9 uses the same banding style as row 7 (bottom border)
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4104) | Out-Null

# Row 10 uses the same banding style as row 6 (no border)
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A10:E10").PasteSpecial(-4104) | Out-Null

# Row 11 has no cell in column A, so only copy columns B:E from row 6
$ws.Range("B6:E6").Copy() | Out-Null
$ws.Range("B11:E11").PasteSpecial(-4104) | Out-Null

$excel.CutCopyMode = 0

# Restore the bottom border on row 9 (copy/paste formatting does not carry borders)
$rowNineBottom = $ws.Range("A9:E9").Borders.Item(9)
$rowNineBottom.LineStyle = 1
$rowNineBottom.Weight = 2

# ---- Write the actual cell values ----
$ws.Cells.Item(8,1).Value2 = $A8
$ws.Cells.Item(8,2).Value2 = $B8
$ws.Cells.Item(8,3).Value2 = $C8
$ws.Cells.Item(8,4).Value2 = $D8
$ws.Cells.Item(8,5).Value2 = $E8

$ws.Cells.Item(9,1).Value2 = $A9
$ws.Cells.Item(9,2).Value2 = $B9
$ws.Cells.Item(9,3).Value2 = $C9
$ws.Cells.Item(9,4).Value2 = $D9
$ws.Cells.Item(9,5).Value2 = $E9

$ws.Cells.Item(10,1).Value2 = $A10
$ws.Cells.Item(10,2).Value2 = $B10
$ws.Cells.Item(10,3).Value2 = $C10
$ws.Cells.Item(10,4).Value2 = $D10
$ws.Cells.Item(10,5).Value2 = $E10

$ws.Cells.Item(11,2).Value2 = $B11
$ws.Cells.Item(11,3).Value2 = $C11
$ws.Cells.Item(11,4).Value2 = $D11
$ws.Cells.Item(11,5).Value2 = $E11

# ---- Row heights (wrapped text rows) ----
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 31.8

# ---- Update view/selection to mirror where the user ended up editing ----
$ws.Range("D12").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
